$d = $word.ActiveDocument

# 1) "Basically" -> "Basically," (the proofErr removal around it happens
#    automatically since this run's text changes and Find/Replace will
#    merge through the grammar-error markers).
$null = $d.Content.Find.Execute("Basically", $true, $false, $false, $false, $false, $true, 1, $false, "Basically,", 2)

Write-Output "done"
